# Auto-generated edit script: update Aegis_Profits (leve profit) sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 109.9
$ws.Range("I8").Value = 109.9
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 329.7
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -190.7
$ws.Range("N8").ClearContents()
$ws.Range("H18").Value = 8163.923
$ws.Range("I18").Value = 510.75
$ws.Range("J18").Value = 100002
$ws.Range("K18").Value = 510.75
$ws.Range("L18").Value = 100002
$ws.Range("M18").Value = -226.75
$ws.Range("N18").Value = -100570
$ws.Range("H32").Value = 1090.3334
$ws.Range("I32").Value = 980
$ws.Range("J32").Value = 1112.4
$ws.Range("K32").Value = 980
$ws.Range("L32").Value = 1112.4
$ws.Range("M32").Value = -654
$ws.Range("N32").Value = -1764.4
$ws.Range("H64").Value = 127837.125
$ws.Range("I64").Value = 335332.66
$ws.Range("J64").Value = 3339.8
$ws.Range("K64").Value = 335332.66
$ws.Range("L64").Value = 3339.8
$ws.Range("M64").Value = -335084.66
$ws.Range("N64").Value = -3835.8
$ws.Range("H67").Value = 127837.125
$ws.Range("I67").Value = 335332.66
$ws.Range("J67").Value = 3339.8
$ws.Range("K67").Value = 335332.66
$ws.Range("L67").Value = 3339.8
$ws.Range("M67").Value = -334474.66
$ws.Range("N67").Value = -5055.8
$ws.Range("H74").Value = 4981.5
$ws.Range("I74").Value = 4957.8
$ws.Range("K74").Value = 4957.8
$ws.Range("M74").Value = -4021.8
$ws.Range("H76").Value = 4078.5881
$ws.Range("I76").Value = 3316.625
$ws.Range("K76").Value = 3316.625
$ws.Range("M76").Value = -3001.625
$ws.Range("H77").Value = 4981.5
$ws.Range("I77").Value = 4957.8
$ws.Range("K77").Value = 24789
$ws.Range("M77").Value = -20109
$ws.Range("H79").Value = 4078.5881
$ws.Range("I79").Value = 3316.625
$ws.Range("K79").Value = 3316.625
$ws.Range("M79").Value = -2224.625
$ws.Range("H107").Value = 284.1154
$ws.Range("I107").Value = 242.71428
$ws.Range("K107").Value = 242.71428
$ws.Range("M107").Value = 1677.28572
$ws.Range("H129").Value = 3699.5557
$ws.Range("J129").Value = 948.1142599999999
$ws.Range("L129").Value = 2844.34278
$ws.Range("N129").Value = -12844.34278

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 12505888
$ws.Range("I6").Value = 20005600
$ws.Range("J6").Value = 6367
$ws.Range("K6").Value = 20005600
$ws.Range("L6").Value = 6367
$ws.Range("M6").Value = -20005427
$ws.Range("N6").Value = -6713
$ws.Range("H63").Value = 2693.4614
$ws.Range("I63").Value = 2164.375
$ws.Range("J63").Value = 3540
$ws.Range("K63").Value = 2164.375
$ws.Range("L63").Value = 3540
$ws.Range("M63").Value = -1478.375
$ws.Range("N63").Value = -4912
$ws.Range("H66").Value = 2693.4614
$ws.Range("I66").Value = 2164.375
$ws.Range("J66").Value = 3540
$ws.Range("K66").Value = 10821.875
$ws.Range("L66").Value = 17700
$ws.Range("M66").Value = -7389.875
$ws.Range("N66").Value = -24564

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 126807.31
$ws.Range("I105").Value = 92809.91
$ws.Range("K105").Value = 92809.91
$ws.Range("M105").Value = -91062.91

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2670.4285
$ws.Range("I62").Value = 2598
$ws.Range("K62").Value = 2598
$ws.Range("M62").Value = -1974
$ws.Range("H65").Value = 2670.4285
$ws.Range("I65").Value = 2598
$ws.Range("K65").Value = 12990
$ws.Range("M65").Value = -9870
$ws.Range("H86").Value = 3399.5
$ws.Range("I86").Value = 3000
$ws.Range("K86").Value = 3000
$ws.Range("M86").Value = -1877
$ws.Range("H89").Value = 3399.5
$ws.Range("I89").Value = 3000
$ws.Range("K89").Value = 15000
$ws.Range("M89").Value = -9384

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 38.47059
$ws.Range("J12").Value = 40.8125
$ws.Range("L12").Value = 122.4375
$ws.Range("N12").Value = -468.4375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 56604.28
$ws.Range("J70").Value = 5748.2
$ws.Range("L70").Value = 5748.2
$ws.Range("N70").Value = -6288.2
$ws.Range("H73").Value = 56604.28
$ws.Range("J73").Value = 5748.2
$ws.Range("L73").Value = 5748.2
$ws.Range("N73").Value = -7620.2
$ws.Range("H80").Value = 125003610
$ws.Range("I80").Value = 250003710
$ws.Range("J80").Value = 3507.5
$ws.Range("K80").Value = 250003710
$ws.Range("L80").Value = 3507.5
$ws.Range("M80").Value = -250002712
$ws.Range("N80").Value = -5503.5
$ws.Range("H83").Value = 125003610
$ws.Range("I83").Value = 250003710
$ws.Range("J83").Value = 3507.5
$ws.Range("K83").Value = 1250018550
$ws.Range("L83").Value = 17537.5
$ws.Range("M83").Value = -1250013558
$ws.Range("N83").Value = -27521.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 15999.556
$ws.Range("J74").Value = 16124.5
$ws.Range("L74").Value = 16124.5
$ws.Range("N74").Value = -18120.5
$ws.Range("H75").Value = 20000
$ws.Range("J75").Value = 20000
$ws.Range("L75").Value = 20000
$ws.Range("N75").Value = -21872
$ws.Range("H77").Value = 15999.556
$ws.Range("J77").Value = 16124.5
$ws.Range("L77").Value = 48373.5
$ws.Range("N77").Value = -58357.5
$ws.Range("H78").Value = 20000
$ws.Range("J78").Value = 20000
$ws.Range("L78").Value = 60000
$ws.Range("N78").Value = -69360
$ws.Range("H80").Value = 9976.866
$ws.Range("J80").Value = 9976.866
$ws.Range("L80").Value = 9976.866
$ws.Range("N80").Value = -12222.866
$ws.Range("H83").Value = 9976.866
$ws.Range("J83").Value = 9976.866
$ws.Range("L83").Value = 29930.598
$ws.Range("N83").Value = -41162.598

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 82714.5
$ws.Range("J46").Value = 82714.5
$ws.Range("L46").Value = 82714.5
$ws.Range("N46").Value = -83176.5
$ws.Range("H121").Value = 39223.555
$ws.Range("J121").Value = 39223.555
$ws.Range("L121").Value = 39223.555
$ws.Range("N121").Value = -42717.555
$ws.Range("H134").Value = 82714.5
$ws.Range("J134").Value = 82714.5
$ws.Range("L134").Value = 248143.5
$ws.Range("N134").Value = -253213.5

Write-Output "Applied all Aegis_Profits updates"